$d = $word.ActiveDocument

# Map of old text -> new text for this update.
$replacements = [ordered]@{
    "2026-02-21 Saturday" = "2026-02-22 Sunday"
    "43÷7="               = "39÷9="
    "38÷2="               = "46÷2="
    "13÷9="               = "31÷5="
    "53÷8="               = "88÷7="
    "64÷8="               = "25÷3="
    "59÷7="               = "71÷4="
    "28÷9="               = "26÷7="
    "65÷9="               = "16÷3="
    "20÷7="               = "21÷9="
    "11÷6="               = "35÷2="
    "63÷5="               = "52÷4="
    "65÷3="               = "64÷9="
    "38÷9="               = "44÷2="
    "58÷7="               = "74÷7="
    "68÷2="               = "64÷9="
    "95÷7="               = "16÷5="
    "25÷5="               = "34÷8="
    "94÷8="               = "32÷4="
    "81÷8="               = "84÷3="
    "76÷9="               = "36÷5="
    "41÷4="               = "22÷7="
    "24÷7="               = "63÷4="
    "30÷7="               = "64÷2="
    "52÷2="               = "92÷2="
    "23÷2="               = "94÷9="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
